$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 388.41177
$ws.Cells.Item(38, 9).Value = 131.33333
$ws.Cells.Item(38, 10).Value = 1005.4
$ws.Cells.Item(38, 11).Value = 393.99999
$ws.Cells.Item(38, 12).Value = 3016.2
$ws.Cells.Item(38, 13).Value = -21.99998999999997
$ws.Cells.Item(38, 14).Value = -3760.2

$ws.Cells.Item(58, 8).Value = 1149.9
$ws.Cells.Item(58, 9).Value = 93
$ws.Cells.Item(58, 10).Value = 2206.8
$ws.Cells.Item(58, 11).Value = 279
$ws.Cells.Item(58, 12).Value = 6620.400000000001
$ws.Cells.Item(58, 13).Value = -129
$ws.Cells.Item(58, 14).Value = -6920.400000000001

$ws.Cells.Item(118, 8).Value = 905.41174
$ws.Cells.Item(118, 9).Value = 480
$ws.Cells.Item(118, 10).Value = 1384
$ws.Cells.Item(118, 11).Value = 1440
$ws.Cells.Item(118, 12).Value = 4152
$ws.Cells.Item(118, 13).Value = 217
$ws.Cells.Item(118, 14).Value = -7466

$ws.Cells.Item(137, 8).Value = 6669456.5
$ws.Cells.Item(137, 9).Value = 10002110
$ws.Cells.Item(137, 10).Value = 4150
$ws.Cells.Item(137, 11).Value = 30006330
$ws.Cells.Item(137, 12).Value = 12450
$ws.Cells.Item(137, 13).Value = -30003780
$ws.Cells.Item(137, 14).Value = -17550

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 326.375
$ws.Cells.Item(22, 9).Value = 333.33334
$ws.Cells.Item(22, 11).Value = 333.33334
$ws.Cells.Item(22, 13).Value = -160.33334

$ws.Cells.Item(64, 8).Value = 1071.8334
$ws.Cells.Item(64, 9).Value = 500
$ws.Cells.Item(64, 10).Value = 1186.2
$ws.Cells.Item(64, 11).Value = 500
$ws.Cells.Item(64, 12).Value = 1186.2
$ws.Cells.Item(64, 13).Value = -275
$ws.Cells.Item(64, 14).Value = -1636.2

$ws.Cells.Item(67, 8).Value = 1071.8334
$ws.Cells.Item(67, 9).Value = 500
$ws.Cells.Item(67, 10).Value = 1186.2
$ws.Cells.Item(67, 11).Value = 500
$ws.Cells.Item(67, 12).Value = 1186.2
$ws.Cells.Item(67, 13).Value = 280
$ws.Cells.Item(67, 14).Value = -2746.2

$ws.Cells.Item(134, 8).Value = 1176.5454
$ws.Cells.Item(134, 9).Value = 1056.3158
$ws.Cells.Item(134, 10).Value = 1938
$ws.Cells.Item(134, 11).Value = 3168.9474
$ws.Cells.Item(134, 12).Value = 5814
$ws.Cells.Item(134, 13).Value = -633.9474
$ws.Cells.Item(134, 14).Value = -10884

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 754.2857
$ws.Cells.Item(22, 9).Value = 293.33334
$ws.Cells.Item(22, 10).Value = 1100
$ws.Cells.Item(22, 11).Value = 293.33334
$ws.Cells.Item(22, 12).Value = 1100
$ws.Cells.Item(22, 13).Value = 56.66665999999998
$ws.Cells.Item(22, 14).Value = -1800

$ws.Cells.Item(31, 8).Value = 1927.8182
$ws.Cells.Item(31, 9).Value = 1505.5333
$ws.Cells.Item(31, 10).Value = 2832.7144
$ws.Cells.Item(31, 11).Value = 1505.5333
$ws.Cells.Item(31, 12).Value = 2832.7144
$ws.Cells.Item(31, 13).Value = -1210.5333
$ws.Cells.Item(31, 14).Value = -3422.7144

$ws.Cells.Item(34, 8).Value = 1927.8182
$ws.Cells.Item(34, 9).Value = 1505.5333
$ws.Cells.Item(34, 10).Value = 2832.7144
$ws.Cells.Item(34, 11).Value = 1505.5333
$ws.Cells.Item(34, 12).Value = 2832.7144
$ws.Cells.Item(34, 13).Value = -1303.5333
$ws.Cells.Item(34, 14).Value = -3236.7144

$ws.Cells.Item(58, 8).Value = 27779576
$ws.Cells.Item(58, 9).Value = 38462630
$ws.Cells.Item(58, 11).Value = 38462630
$ws.Cells.Item(58, 13).Value = -38462427

$ws.Cells.Item(132, 8).Value = 64818.438
$ws.Cells.Item(132, 9).Value = 2110.8572
$ws.Cells.Item(132, 10).Value = 113591
$ws.Cells.Item(132, 11).Value = 6332.571599999999
$ws.Cells.Item(132, 12).Value = 340773
$ws.Cells.Item(132, 13).Value = -3802.571599999999
$ws.Cells.Item(132, 14).Value = -345833

$ws.Cells.Item(136, 8).Value = 27779576
$ws.Cells.Item(136, 9).Value = 38462630
$ws.Cells.Item(136, 11).Value = 115387890
$ws.Cells.Item(136, 13).Value = -115385340

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1067.3334
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1067.3334
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 3202.0002
$ws.Cells.Item(113, 14).Value = -7542.0002
$ws.Cells.Item(113, 13).ClearContents()

$ws.Cells.Item(131, 8).Value = 1254.72
$ws.Cells.Item(131, 9).Value = 768.75
$ws.Cells.Item(131, 10).Value = 1483.4117
$ws.Cells.Item(131, 11).Value = 2306.25
$ws.Cells.Item(131, 12).Value = 4450.2351
$ws.Cells.Item(131, 13).Value = 2733.75
$ws.Cells.Item(131, 14).Value = -14530.2351

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 4500
$ws.Cells.Item(43, 9).Value = 1625
$ws.Cells.Item(43, 10).Value = 16000
$ws.Cells.Item(43, 11).Value = 1625
$ws.Cells.Item(43, 12).Value = 16000
$ws.Cells.Item(43, 13).Value = -1474
$ws.Cells.Item(43, 14).Value = -16302

$ws.Cells.Item(46, 8).Value = 24762.363
$ws.Cells.Item(46, 10).Value = 24762.363
$ws.Cells.Item(46, 12).Value = 24762.363
$ws.Cells.Item(46, 14).Value = -25074.363

$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 13).ClearContents()

$ws.Cells.Item(62, 8).Value = 24333.334
$ws.Cells.Item(62, 10).Value = 24333.334
$ws.Cells.Item(62, 12).Value = 24333.334
$ws.Cells.Item(62, 14).Value = -25705.334

$ws.Cells.Item(63, 8).Value = 20106.875
$ws.Cells.Item(63, 10).Value = 20106.875
$ws.Cells.Item(63, 12).Value = 20106.875
$ws.Cells.Item(63, 14).Value = -21478.875

$ws.Cells.Item(65, 8).Value = 24333.334
$ws.Cells.Item(65, 10).Value = 24333.334
$ws.Cells.Item(65, 12).Value = 73000.00199999999
$ws.Cells.Item(65, 14).Value = -79864.00199999999

$ws.Cells.Item(66, 8).Value = 20106.875
$ws.Cells.Item(66, 10).Value = 20106.875
$ws.Cells.Item(66, 12).Value = 60320.625
$ws.Cells.Item(66, 14).Value = -67184.625

$ws.Cells.Item(74, 8).Value = 23886.666
$ws.Cells.Item(74, 10).Value = 23886.666
$ws.Cells.Item(74, 12).Value = 23886.666
$ws.Cells.Item(74, 14).Value = -25758.666

$ws.Cells.Item(75, 8).Value = 36118.715
$ws.Cells.Item(75, 10).Value = 36118.715
$ws.Cells.Item(75, 12).Value = 36118.715
$ws.Cells.Item(75, 14).Value = -37866.715

$ws.Cells.Item(77, 8).Value = 23886.666
$ws.Cells.Item(77, 10).Value = 23886.666
$ws.Cells.Item(77, 12).Value = 71659.99800000001
$ws.Cells.Item(77, 14).Value = -81019.99800000001

$ws.Cells.Item(78, 8).Value = 36118.715
$ws.Cells.Item(78, 10).Value = 36118.715
$ws.Cells.Item(78, 12).Value = 108356.145
$ws.Cells.Item(78, 14).Value = -117092.145

$ws.Cells.Item(82, 8).Value = 38000
$ws.Cells.Item(82, 10).Value = 38000
$ws.Cells.Item(82, 12).Value = 38000
$ws.Cells.Item(82, 14).Value = -38766

$ws.Cells.Item(85, 8).Value = 38000
$ws.Cells.Item(85, 10).Value = 38000
$ws.Cells.Item(85, 12).Value = 38000
$ws.Cells.Item(85, 14).Value = -40652

$ws.Cells.Item(86, 8).Value = 37223.145
$ws.Cells.Item(86, 10).Value = 37223.145
$ws.Cells.Item(86, 12).Value = 37223.145
$ws.Cells.Item(86, 14).Value = -39595.145

$ws.Cells.Item(87, 8).Value = 61600
$ws.Cells.Item(87, 10).Value = 61600
$ws.Cells.Item(87, 12).Value = 61600
$ws.Cells.Item(87, 14).Value = -64096

$ws.Cells.Item(88, 8).Value = 40011.43
$ws.Cells.Item(88, 10).Value = 40011.43
$ws.Cells.Item(88, 12).Value = 40011.43
$ws.Cells.Item(88, 14).Value = -40913.43

$ws.Cells.Item(89, 8).Value = 37223.145
$ws.Cells.Item(89, 10).Value = 37223.145
$ws.Cells.Item(89, 12).Value = 111669.435
$ws.Cells.Item(89, 14).Value = -123525.435

$ws.Cells.Item(90, 8).Value = 61600
$ws.Cells.Item(90, 10).Value = 61600
$ws.Cells.Item(90, 12).Value = 184800
$ws.Cells.Item(90, 14).Value = -197280

$ws.Cells.Item(91, 8).Value = 40011.43
$ws.Cells.Item(91, 10).Value = 40011.43
$ws.Cells.Item(91, 12).Value = 40011.43
$ws.Cells.Item(91, 14).Value = -43131.43

$ws.Cells.Item(126, 8).Value = 1715.3182
$ws.Cells.Item(126, 9).Value = 1455.6154
$ws.Cells.Item(126, 11).Value = 4366.8462
$ws.Cells.Item(126, 13).Value = -1896.8462

$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1964.2142
$ws.Cells.Item(61, 9).Value = 2007.6154
$ws.Cells.Item(61, 10).Value = 1400
$ws.Cells.Item(61, 11).Value = 2007.6154
$ws.Cells.Item(61, 12).Value = 1400
$ws.Cells.Item(61, 13).Value = -1805.6154
$ws.Cells.Item(61, 14).Value = -1804

$ws.Cells.Item(113, 8).Value = 1964.2142
$ws.Cells.Item(113, 9).Value = 2007.6154
$ws.Cells.Item(113, 10).Value = 1400
$ws.Cells.Item(113, 11).Value = 2007.6154
$ws.Cells.Item(113, 12).Value = 1400
$ws.Cells.Item(113, 13).Value = 162.3846000000001
$ws.Cells.Item(113, 14).Value = -5740

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 13663.125
$ws.Cells.Item(54, 10).Value = 13663.125
$ws.Cells.Item(54, 12).Value = 13663.125
$ws.Cells.Item(54, 14).Value = -14703.125
